# Add two new columns (I: I0, J: IF) to the worksheet, mirroring the
# existing style/formatting used by the other header cells and filling in
# the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1), using same style as the existing header cells (H1 "IP")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for rows 2-23 : column I (I0) and column J (IF)
$data = @{
    2  = @(1, 4)
    3  = @(7, 8)
    4  = @(1, 2)
    5  = @(1, 3)
    6  = @(1, 5)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 4)
    10 = @(1, 5)
    11 = @(1, 4)
    12 = @(1, 4)
    13 = @(1, 7)
    14 = @(1, 6)
    15 = @(1, 6)
    16 = @(1, 6)
    17 = @(1, 5)
    18 = @(1, 3)
    19 = @(1, 4)
    20 = @(1, 5)
    21 = @(1, 5)
    22 = @(1, 4)
    23 = @(1, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
